$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Solder connections for wifi" columns (E:H) added next to the
# existing ATWINC1500 WiFi block (rows 30-34)

# Row 30 - header row for the taped-bundle columns
$ws.Range("E30").Value = "Yellow Taped Bundle"
$ws.Range("G30").Value = "Non Taped"

# Row 31 - GND / SCK
$ws.Range("E31").Value = "Red"
$ws.Range("F31").Value = "5v and EN"
$ws.Range("G31").Value = "Red"
$ws.Range("H31").Value = "SCK"

# Row 32 - SCK / MISO
$ws.Range("E32").Value = "Black"
$ws.Range("F32").Value = "GND"
$ws.Range("G32").Value = "Black"
$ws.Range("H32").Value = "MISO"

# Row 33 - MISO / MOSI
$ws.Range("E33").Value = "Yellow"
$ws.Range("F33").Value = "IRQ"
$ws.Range("G33").Value = "Yellow"
$ws.Range("H33").Value = "MOSI"

# Row 34 - MOSI / CS
$ws.Range("E34").Value = "White"
$ws.Range("F34").Value = "RST"
$ws.Range("G34").Value = "White"
$ws.Range("H34").Value = "CS"

# Resize columns E and F to fit the newly added content
$ws.Columns.Item(5).ColumnWidth = 19.21875
$ws.Columns.Item(6).ColumnWidth = 11.5546875

# Update the view: selected cells (mirrors the author scrolling to the
# new wifi-solder block and selecting the new H column entries)
$ws.Range("H31:H34").Select()
